$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (A1:I1) ---
# Insert the new "Status" header at D1 and shift month headers; add new "Oct_2025" at G1
# and re-point I1 to "QoQ". Write the header text values first, then copy the formatting
# from an existing styled header cell (A1) onto the two brand-new header cells (D1, I1)
# so they pick up the same bold/centered/bordered style used by the rest of row 1.
$ws.Cells.Item(1, 1).Value = "ISIN"
$ws.Cells.Item(1, 2).Value = "Stock Name"
$ws.Cells.Item(1, 3).Value = "Mutual Fund"
$ws.Cells.Item(1, 4).Value = "Status"
$ws.Cells.Item(1, 5).Value = "Jan_2026"
$ws.Cells.Item(1, 6).Value = "Dec_2025"
$ws.Cells.Item(1, 7).Value = "Oct_2025"
$ws.Cells.Item(1, 8).Value = "MoM"
$ws.Cells.Item(1, 9).Value = "QoQ"

$ws.Range("A1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("I1").PasteSpecial(-4122)

# --- Data rows (A2:I66) ---
# Full refreshed holdings data from the quant engine: a new "Status" column (D) was
# added, the trailing month column rolled from Nov_2025 to Oct_2025 (G), and MoM/QoQ
# (H/I) were recomputed; several positions were fully exited (rows pushed to the bottom).
$data = @(
  @("INE271C01023", "DLF Limited", "quant Arbitrage Fund", "Fresh Entry", 5.122585, 0, 2.18937, 5.122585, 2.933215),
  @("INE742F01042", "Adani Ports & Special Economic Zone Ltd", "quant Arbitrage Fund", "Reducing Consistently", 4.586708, 6.368688, 6.935659, -1.78198, -2.348951),
  @("INE619A01035", "Patanjali Foods Limited", "quant Arbitrage Fund", "Reducing Consistently", 4.469772, 6.535528, 5.580702, -2.065756, -1.11093),
  @("INE296A01032", "Bajaj Finance Limited", "quant Arbitrage Fund", "Reducing Consistently", 4.246653, 5.613014, 4.894347, -1.366360999999999, -0.6476939999999995),
  @("INE752E01010", "Power Grid Corporation of India Limited", "quant Arbitrage Fund", "Fresh Entry", 3.873369, 0, 0, 3.873369, 3.873369),
  @("INE669E01016", "Vodafone Idea Ltd.", "quant Arbitrage Fund", "Reducing Consistently", 3.851404, 4.976178, 4.452223, -1.124774, -0.600819),
  @("INE038A01020", "Hindalco Industries Limited", "quant Arbitrage Fund", "Reducing", 3.783416, 4.674474, 2.290919, -0.8910580000000001, 1.492497),
  @("INE205A01025", "Vedanta Limited", "quant Arbitrage Fund", "Reducing", 3.130176, 3.723174, 1.991723, -0.5929980000000001, 1.138453),
  @("INE115A01026", "LIC Housing Finance Ltd", "quant Arbitrage Fund", "Adding Consistently", 2.97592, 2.060316, 2.404248, 0.9156040000000001, 0.571672),
  @("INE406A01037", "Aurobindo Pharma Limited", "quant Arbitrage Fund", "Reducing Consistently", 2.363829, 3.105701, 3.480331, -0.7418719999999999, -1.116502),
  @("INE154A01025", "ITC Limited", "quant Arbitrage Fund", "Reducing", 2.241933, 3.18037, 0, -0.938437, 2.241933),
  @("INE016A01026", "Dabur India Limited", "quant Arbitrage Fund", "Fresh Entry", 2.003167, 0, 0.285146, 2.003167, 1.718021),
  @("INE776C01039", "GMR Airports Limited", "quant Arbitrage Fund", "Adding Consistently", 1.968018, 1.042445, 1.571687, 0.925573, 0.396331),
  @("INE094A01015", "Hindustan Petroleum Corporation Ltd", "quant Arbitrage Fund", "Reducing", 1.948855, 3.055017, 0.39462, -1.106162, 1.554235),
  @("INE918I01026", "Bajaj Finserv Ltd.", "quant Arbitrage Fund", "Fresh Entry", 1.930213, 0, 0, 1.930213, 1.930213),
  @("INE245A01021", "Tata Power Company Limited", "quant Arbitrage Fund", "Adding Consistently", 1.84816, 1.576368, 0, 0.271792, 1.84816),
  @("INE019A01038", "JSW Steel Limited", "quant Arbitrage Fund", "Reducing Consistently", 1.847533, 2.376838, 2.71378, -0.5293049999999997, -0.8662469999999998),
  @("INE002A01018", "Reliance Industries Limited", "quant Arbitrage Fund", "Reducing Consistently", 1.682869, 2.54027, 2.651455, -0.8574010000000001, -0.968586),
  @("INE216A01030", "Britannia Industries Limited", "quant Arbitrage Fund", "Fresh Entry", 1.651093, 0, 0, 1.651093, 1.651093),
  @("INE849A01020", "Trent Ltd", "quant Arbitrage Fund", "Reducing", 1.496839, 2.269405, 0, -0.7725659999999999, 1.496839),
  @("INE237A01036", "Kotak Mahindra Bank Limited", "quant Arbitrage Fund", "Fresh Entry", 1.45196, 0, 0, 1.45196, 1.45196),
  @("INE030A01027", "Hindustan Unilever Limited", "quant Arbitrage Fund", "Fresh Entry", 1.435625, 0, 0, 1.435625, 1.435625),
  @("INE121J01017", "Indus Towers Limited", "quant Arbitrage Fund", "Reducing", 1.403702, 1.774481, 0, -0.370779, 1.403702),
  @("INE090A01021", "ICICI Bank Limited", "quant Arbitrage Fund", "Reducing Consistently", 1.312675, 1.744936, 1.927675, -0.432261, -0.615),
  @("INE361B01024", "Divi's Laboratories Limited", "quant Arbitrage Fund", "Reducing Consistently", 1.196325, 1.695159, 1.970379, -0.498834, -0.774054),
  @("INE0J1Y01017", "Life Insurance Corporation Of India", "quant Arbitrage Fund", "Fresh Entry", 1.163887, 0, 0, 1.163887, 1.163887),
  @("INE745G01043", "Multi Commodity Exchange of India Ltd.", "quant Arbitrage Fund", "Fresh Entry", 1.124557, 0, 0, 1.124557, 1.124557),
  @("INE040A01034", "HDFC Bank Limited", "quant Arbitrage Fund", "Reducing Consistently", 1.091292, 1.561304, 1.969029, -0.4700120000000001, -0.877737),
  @("INE257A01026", "Bharat Heavy Electricals Ltd", "quant Arbitrage Fund", "Fresh Entry", 1.09069, 0, 0, 1.09069, 1.09069),
  @("INE118H01025", "BSE Ltd", "quant Arbitrage Fund", "Reducing", 1.036849, 1.308759, 0, -0.2719100000000001, 1.036849),
  @("INE028A01039", "Bank of Baroda", "quant Arbitrage Fund", "Reducing", 0.969589, 1.285282, 0, -0.315693, 0.969589),
  @("INE263A01024", "Bharat Electronics Ltd", "quant Arbitrage Fund", "Reducing Consistently", 0.860185, 1.026807, 1.207408, -0.166622, -0.3472230000000001),
  @("INE376G01013", "Biocon Ltd", "quant Arbitrage Fund", "Fresh Entry", 0.725189, 0, 0, 0.725189, 0.725189),
  @("INE303R01014", "Kalyan Jewellers India Limited", "quant Arbitrage Fund", "Reducing", 0.722515, 1.300562, 0, -0.578047, 0.722515),
  @("INE047A01021", "Grasim Industries Ltd", "quant Arbitrage Fund", "Fresh Entry", 0.69667, 0, 0, 0.69667, 0.69667),
  @("INE584A01023", "NMDC Ltd", "quant Arbitrage Fund", "Reducing", 0.650259, 0.893226, 0.299202, -0.2429669999999999, 0.351057),
  @("INE067A01029", "CG Power and Industrial Solutions Ltd", "quant Arbitrage Fund", "Reducing Consistently", 0.647791, 0.963852, 1.794423, -0.316061, -1.146632),
  @("INE062A01020", "State Bank of India", "quant Arbitrage Fund", "Reducing", 0.574992, 0.70324, 0.411007, -0.128248, 0.1639849999999999),
  @("INE298A01020", "Cummins India Ltd.", "quant Arbitrage Fund", "Reducing", 0.552887, 0.79962, 0, -0.246733, 0.552887),
  @("INE881D01027", "Oracle Financial Services Software Ltd", "quant Arbitrage Fund", "Reducing Consistently", 0.552314, 0.733836, 0.896408, -0.1815220000000001, -0.344094),
  @("INE095N01031", "National Building Construction Corp", "quant Arbitrage Fund", "Reducing", 0.508178, 0.839563, 0, -0.3313849999999999, 0.508178),
  @("INE192A01025", "Tata Consumer Products Ltd", "quant Arbitrage Fund", "Fresh Entry", 0.369897, 0, 0, 0.369897, 0.369897),
  @("INE160A01022", "Punjab National Bank", "quant Arbitrage Fund", "Fresh Entry", 0.356413, 0, 0, 0.356413, 0.356413),
  @("INE031A01017", "Housing & Urban Devlopment Company Ltd", "quant Arbitrage Fund", "Reducing Consistently", 0.336468, 0.537199, 0.615432, -0.200731, -0.278964),
  @("INE148O01028", "Delhivery Limited", "quant Arbitrage Fund", "Reducing", 0.312321, 0.399991, 0, -0.08766999999999997, 0.312321),
  @("INE484J01027", "Godrej Properties Limited", "quant Arbitrage Fund", "Fresh Entry", 0.308627, 0, 0.625585, 0.308627, -0.316958),
  @("INE467B01029", "Tata Consultancy Services Limited", "quant Arbitrage Fund", "Fresh Entry", 0.281016, 0, 0, 0.281016, 0.281016),
  @("INE947Q01028", "Laurus Labs Ltd", "quant Arbitrage Fund", "Fresh Entry", 0.162169, 0, 1.896342, 0.162169, -1.734173),
  @("INE476A01022", "Canara Bank", "quant Arbitrage Fund", "Complete Exit", 0, 0, 1.514256, 0, -1.514256),
  @("INE134E01011", "Power Finance Corporation Ltd.", "quant Arbitrage Fund", "Complete Exit", 0, 0, 2.360789, 0, -2.360789),
  @("INE029A01011", "Bharat Petroleum Corp Ltd", "quant Arbitrage Fund", "Complete Exit", 0, 0.281557, 1.813399, -0.281557, -1.813399),
  @("INE059A01026", "Cipla Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0.420804, 0, -0.420804, 0),
  @("INE066F01020", "Hindustan Aeronautics Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0, 0.615827, 0, -0.615827),
  @("INE081A01020", "Tata Steel Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0, 1.882056, 0, -1.882056),
  @("INE758E01017", "Jio Financial Services Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0, 2.74085, 0, -2.74085),
  @("INE129A01019", "GAIL (India) Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0.92037, 0, -0.92037, 0),
  @("INE745G01035", "Multi Commodity Exchange of India Ltd.", "quant Arbitrage Fund", "Complete Exit", 0, 1.328867, 0, -1.328867, 0),
  @("INE397D01024", "Bharti Airtel Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0, 0.057075, 0, -0.057075),
  @("INE148I01020", "Sammaan Capital Ltd.", "quant Arbitrage Fund", "Complete Exit", 0, 0, 2.181693, 0, -2.181693),
  @("INE721A01047", "Shriram Finance Limited", "quant Arbitrage Fund", "Complete Exit", 0, 1.351238, 3.468947, -1.351238, -3.468947),
  @("INE010B01027", "Zydus Lifesciences Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0, 0.615506, 0, -0.615506),
  @("INE280A01028", "Titan Company Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0, 0.345126, 0, -0.345126),
  @("INE596I01020", "Computer Age Management Services Ltd", "quant Arbitrage Fund", "Complete Exit", 0, 0.530473, 0, -0.530473, 0),
  @("INE326A01037", "Lupin Limited", "quant Arbitrage Fund", "Complete Exit", 0, 0, 0.976109, 0, -0.976109),
  @("INE237A01028", "Kotak Mahindra Bank Limited", "quant Arbitrage Fund", "Complete Exit", 0, 2.10127, 0, -2.10127, 0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $excelRow = $i + 2
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($excelRow, $j + 1).Value = $row[$j]
    }
}